$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1033.5
$ws.Range("I17").Value = 554.1429000000001
$ws.Range("J17").Value = 1369.05
$ws.Range("K17").Value = 1662.4287
$ws.Range("L17").Value = 4107.15
$ws.Range("M17").Value = -1494.4287
$ws.Range("N17").Value = -4443.15

$ws.Range("H62").Value = 5103.7915
$ws.Range("I62").Value = 2373.75
$ws.Range("J62").Value = 7833.8335
$ws.Range("K62").Value = 2373.75
$ws.Range("L62").Value = 7833.8335
$ws.Range("M62").Value = -1749.75
$ws.Range("N62").Value = -9081.833500000001

$ws.Range("H65").Value = 5103.7915
$ws.Range("I65").Value = 2373.75
$ws.Range("J65").Value = 7833.8335
$ws.Range("K65").Value = 11868.75
$ws.Range("L65").Value = 39169.1675
$ws.Range("M65").Value = -8748.75
$ws.Range("N65").Value = -45409.1675

$ws.Range("H116").Value = 2652.158
$ws.Range("I116").Value = 2498.9285
$ws.Range("K116").Value = 2498.9285
$ws.Range("M116").Value = 943.0715

$ws.Range("H132").Value = 3976.8948
$ws.Range("I132").Value = 3904.1333
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 11712.3999
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -9182.3999
$ws.Range("N132").Value = -17809.25

$ws.Range("H137").Value = 1295.5
$ws.Range("I137").Value = 1023.125
$ws.Range("K137").Value = 3069.375
$ws.Range("M137").Value = -519.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 18184.2
$ws.Range("I31").Value = 1691.7142
$ws.Range("K31").Value = 1691.7142
$ws.Range("M31").Value = -1397.7142

$ws.Range("H32").Value = 1776426.6
$ws.Range("I32").Value = 2030524.9
$ws.Range("K32").Value = 2030524.9
$ws.Range("M32").Value = -2030237.9

$ws.Range("H63").Value = 30421.295
$ws.Range("I63").Value = 67463.42999999999
$ws.Range("J63").Value = 4491.8
$ws.Range("K63").Value = 67463.42999999999
$ws.Range("L63").Value = 4491.8
$ws.Range("M63").Value = -66777.42999999999
$ws.Range("N63").Value = -5863.8

$ws.Range("H66").Value = 30421.295
$ws.Range("I66").Value = 67463.42999999999
$ws.Range("J66").Value = 4491.8
$ws.Range("K66").Value = 337317.15
$ws.Range("L66").Value = 22459
$ws.Range("M66").Value = -333885.15
$ws.Range("N66").Value = -29323

$ws.Range("H74").Value = 2317.776
$ws.Range("I74").Value = 1257.697
$ws.Range("J74").Value = 3717.08
$ws.Range("K74").Value = 1257.697
$ws.Range("L74").Value = 3717.08
$ws.Range("M74").Value = -383.6969999999999
$ws.Range("N74").Value = -5465.08

$ws.Range("H77").Value = 2317.776
$ws.Range("I77").Value = 1257.697
$ws.Range("J77").Value = 3717.08
$ws.Range("K77").Value = 6288.485
$ws.Range("L77").Value = 18585.4
$ws.Range("M77").Value = -1920.485
$ws.Range("N77").Value = -27321.4

$ws.Range("H80").Value = 17999.25
$ws.Range("J80").Value = 17999.25
$ws.Range("L80").Value = 17999.25
$ws.Range("N80").Value = -19995.25

$ws.Range("H83").Value = 17999.25
$ws.Range("J83").Value = 17999.25
$ws.Range("L83").Value = 53997.75
$ws.Range("N83").Value = -63981.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41668570
$ws.Range("I105").Value = 62501404
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 62501404
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -62499657
$ws.Range("N105").Value = -6394

$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4123.4385
$ws.Range("I31").Value = 1092.8334
$ws.Range("J31").Value = 7072.1353
$ws.Range("K31").Value = 1092.8334
$ws.Range("L31").Value = 7072.1353
$ws.Range("M31").Value = -797.8334
$ws.Range("N31").Value = -7662.1353

$ws.Range("H34").Value = 4123.4385
$ws.Range("I34").Value = 1092.8334
$ws.Range("J34").Value = 7072.1353
$ws.Range("K34").Value = 1092.8334
$ws.Range("L34").Value = 7072.1353
$ws.Range("M34").Value = -890.8334
$ws.Range("N34").Value = -7476.1353

$ws.Range("H86").Value = 3401.1667
$ws.Range("I86").Value = 3541.4
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 3541.4
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -2418.4
$ws.Range("N86").Value = -4946

$ws.Range("H89").Value = 3401.1667
$ws.Range("I89").Value = 3541.4
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 17707
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -12091
$ws.Range("N89").Value = -24732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4011.1667
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 3953
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 11859
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -13731

$ws.Range("H83").Value = 4011.1667
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 3953
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 35577
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -44937

$ws.Range("H131").Value = 3531.451
$ws.Range("I131").Value = 571.6667
$ws.Range("J131").Value = 3926.0889
$ws.Range("K131").Value = 1715.0001
$ws.Range("L131").Value = 11778.2667
$ws.Range("M131").Value = 3324.9999
$ws.Range("N131").Value = -21858.2667

$ws.Range("H137").Value = 30825.025
$ws.Range("J137").Value = 56100.684
$ws.Range("L137").Value = 168302.052
$ws.Range("N137").Value = -178502.052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 38655
$ws.Range("J118").Value = 38655
$ws.Range("L118").Value = 38655
$ws.Range("N118").Value = -41969

$ws.Range("H132").Value = 3547.4285
$ws.Range("I132").Value = 2776.3809
$ws.Range("K132").Value = 8329.1427
$ws.Range("M132").Value = -5799.1427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1471
$ws.Range("I113").Value = 1579.75
$ws.Range("J113").Value = 1210
$ws.Range("K113").Value = 4739.25
$ws.Range("L113").Value = 3630
$ws.Range("M113").Value = -2569.25
$ws.Range("N113").Value = -7970
